$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B to fit the time values and apply the HH:MM:SS time format
# to the whole column (this also re-stamps the existing B1:B4 cells,
# including the header B1).
$ws.Columns.Item(2).ColumnWidth = 38.16
$ws.Columns.Item(2).NumberFormat = "HH:MM:SS"

# New row 5: midnight (00:00:00)
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 0
$ws.Range("B5").NumberFormat = "HH:MM:SS"
$ws.Range("C5").Value = 25569
$ws.Range("C5").NumberFormat = "YYYY\-MM\-DD\ HH:MM:SS"

# New row 6: one second to midnight (23:59:59) - time values now carry
# one-second precision
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 0.999988425925926
$ws.Range("B6").NumberFormat = "HH:MM:SS"
$ws.Range("C6").Value = 25569
$ws.Range("C6").NumberFormat = "YYYY\-MM\-DD\ HH:MM:SS"

$ws.Range("B6").Select()
